$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the time-availability data added to the truck matrix.
$ws.Range("I3").Value = "Disponibilidad de tiempo(int)"

# Populate the new column with the availability value for every truck row (4-21).
for ($r = 4; $r -le 21; $r++) {
    $ws.Cells.Item($r, 9).Value = 11
}

# Leave the selection where the editor last left it.
$ws.Range("K17").Select()
